$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Adição do crawler do Twitter ---
# The "Temáticas" table no longer lists ID_Tweet (G4) and the "Tweets"
# table loses its Link_Tweet field (C5). The "Notícias" table's column
# of field names shifts up by one row (the stray duplicate ID_Candidato
# at E3 is removed) and ends up one row shorter overall.

# Remove the stray ID_Tweet entry from the Temáticas table.
$ws.Range("G4").Clear()

# Remove the Link_Tweet entry from the Tweets table.
$ws.Range("C5").Clear()

# Shift the Notícias field names up: drop the duplicate ID_Candidato,
# and move Titulo_Noticia / Conteudo_Noticia / Link_Noticia /
# Relevancia_Noticia up one row each.
$ws.Range("E3").Value2 = "Titulo_Noticia"
$ws.Range("E4").Value2 = "Conteudo_Noticia"
$ws.Range("E5").Value2 = "Link_Noticia"
$ws.Range("E6").Value2 = "Relevancia_Noticia"

# The old last row (row 7) is no longer used.
$ws.Range("E7").Clear()

# Update the view selection to match the new layout.
$ws.Range("G4").Select()
